$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2026-01-19 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-20 Tuesday", 2)

# The division problems live in a single table; rows 1,5,9,13,17 (1-based)
# contain the five problems per "block", the rows in between are blank
# spacer rows. Addressing cells directly by (row, column) avoids any
# ambiguity caused by duplicate expressions appearing more than once.
$table = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; Text = "76÷5=" },
    @{ Row = 1;  Col = 2; Text = "37÷8=" },
    @{ Row = 1;  Col = 3; Text = "59÷2=" },
    @{ Row = 1;  Col = 4; Text = "79÷6=" },
    @{ Row = 1;  Col = 5; Text = "41÷7=" },

    @{ Row = 5;  Col = 1; Text = "56÷6=" },
    @{ Row = 5;  Col = 2; Text = "16÷3=" },
    @{ Row = 5;  Col = 3; Text = "97÷2=" },
    @{ Row = 5;  Col = 4; Text = "95÷4=" },
    @{ Row = 5;  Col = 5; Text = "34÷2=" },

    @{ Row = 9;  Col = 1; Text = "22÷9=" },
    @{ Row = 9;  Col = 2; Text = "66÷9=" },
    @{ Row = 9;  Col = 3; Text = "76÷3=" },
    @{ Row = 9;  Col = 4; Text = "31÷4=" },
    @{ Row = 9;  Col = 5; Text = "53÷4=" },

    @{ Row = 13; Col = 1; Text = "46÷8=" },
    @{ Row = 13; Col = 2; Text = "35÷7=" },
    @{ Row = 13; Col = 3; Text = "13÷9=" },
    @{ Row = 13; Col = 4; Text = "42÷5=" },
    @{ Row = 13; Col = 5; Text = "40÷9=" },

    @{ Row = 17; Col = 1; Text = "79÷6=" },
    @{ Row = 17; Col = 2; Text = "84÷4=" },
    @{ Row = 17; Col = 3; Text = "66÷7=" },
    @{ Row = 17; Col = 4; Text = "71÷4=" },
    @{ Row = 17; Col = 5; Text = "28÷9=" }
)

foreach ($e in $edits) {
    $cell = $table.Cell($e.Row, $e.Col)
    $r = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so we only
    # replace the visible text of the cell.
    $r.End = $r.End - 1
    $r.Text = $e.Text
}
